# Updated cryptos list (Price / Volume(1h) refresh), matching the
# upstream GitHub Actions scraper commit.
#
# Note: several "Price" values are digit strings that Excel's COM layer
# would otherwise auto-coerce to a number (losing formatting / precision,
# e.g. "7.80" -> 7.7999999999999998, "0.999" -> 0.999 but typed as a
# number instead of text). Those are written with a leading apostrophe
# to force text, then the cell style is reset to "Normal" so no stray
# number-format/style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.388.38'
$ws.Range('E2').Value = '  -2.58%  '
$ws.Range('D3').Value = '1.650.83'
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "'213.55"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('E6').Value = '  -2.11%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = "'24.11"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = "'0.263"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('D11').Value = "'0.0877"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('D12').Value = '1.883.18'
$ws.Range('E12').Value = '  -2.58%  '
$ws.Range('D13').Value = '1.649.24'
$ws.Range('E13').Value = '  -2.48%  '
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').Value = "'0.571"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.58%  '
$ws.Range('D16').Value = "'65.85"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.66%  '
$ws.Range('D17').Value = '27.373.45'
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('D18').Value = "'234.01"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.77%  '
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('E20').Value = '  -3.31%  '
$ws.Range('D21').Value = "'0.999"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('E22').Value = '  -3.18%  '
$ws.Range('D23').Value = "'9.29"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.71%  '
$ws.Range('D24').Value = "'2.01"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').Value = "'146.67"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.64%  '
$ws.Range('D26').Value = "'7.18"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.23%  '
$ws.Range('D27').Value = "'16.06"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.87%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  -2.37%  '
$ws.Range('D30').Value = "'0.0497"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('E32').Value = '  -2.42%  '
$ws.Range('D33').Value = '1.460.46'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('E35').Value = '  -4.13%  '
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').Value = "'0.908"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.71%  '
$ws.Range('D38').Value = "'0.574"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.29%  '
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('D41').Value = "'0.999"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('E42').Value = '  -3.15%  '
$ws.Range('D43').Value = "'65.14"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.02%  '
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = "'0.785"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.68%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.792.20'
$ws.Range('E46').Value = '  -2.44%  '
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('D48').Value = "'88.34"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('E49').Value = '  -4.35%  '
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').Value = "'7.80"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.67%  '
